$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the category header rows and trailing footnote rows (bottom-up so
# row indices of not-yet-deleted rows stay valid).
$rowsToDelete = @(35, 34, 27, 19, 13, 8, 5)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
